# Applies the cryptos-list refresh described by the diff (Sun Jul 30 2023 run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values that are NOT plain decimal numbers (coin names, URLs, the
#     "  +x.xx%  " volume strings, and multi-dot price strings like
#     "29.404.56") can be written straight to .Value: Excel will not
#     reinterpret them as a Number, so they stay stored as text, matching
#     the source inline strings. ---
$ws.Range("D2").Value = "29.403.32"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.877.70"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("E11").Value = "  -3.14%  "
$ws.Range("D12").Value = "1.884.40"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  +4.19%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "29.406.61"
$ws.Range("E19").Value = "  +5.13%  "
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "2.136.57"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").Value = "1.273.96"
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").Value = "2.030.77"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("E51").Value = "  -0.40%  "

# --- Price values that DO look like plain decimals must still be stored as
#     TEXT (every Price cell in the source sheet is a string, never a
#     Number). Assigning such a string straight to .Value makes Excel auto-
#     convert it to a Number, so instead write a formula that evaluates to
#     the literal text, then copy/paste-special that single cell back onto
#     itself as a value only. That collapses the cell to a plain text result
#     with no leftover formula and no number-format/style change (each cell
#     is handled individually - Excel does not reliably map values 1:1 when
#     copy/paste-special is done over a large non-contiguous union range). ---
$ws.Range("D5").Formula = "=`"0.7167`""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = "=`"243.80`""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D8").Formula = "=`"0.07953`""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("D9").Formula = "=`"0.3148`""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("D10").Formula = "=`"24.95`""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("D11").Formula = "=`"0.08126`""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D13").Formula = "=`"95.22`""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("D14").Formula = "=`"5.239`""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("D15").Formula = "=`"0.7071`""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("D16").Formula = "=`"6.415`""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("D17").Formula = "=`"0.000008426`""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("D19").Formula = "=`"253.33`""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("D20").Formula = "=`"13.37`""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D26").Formula = "=`"9.065`""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("D27").Formula = "=`"162.01`""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("D28").Formula = "=`"18.90`""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("D30").Formula = "=`"4.419`""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D31").Formula = "=`"4.307`""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("D33").Formula = "=`"0.05328`""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("D34").Formula = "=`"1.949`""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("D35").Formula = "=`"0.7554`""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("D37").Formula = "=`"2.701`""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("D38").Formula = "=`"0.01894`""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("D40").Formula = "=`"2.765`""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("D41").Formula = "=`"6.395`""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("D42").Formula = "=`"0.9069`""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("D43").Formula = "=`"112.14`""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("D44").Formula = "=`"74.31`""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D48").Formula = "=`"1.805`""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("D49").Formula = "=`"0.5201`""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("D50").Formula = "=`"9.524`""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("D51").Formula = "=`"0.4346`""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

$excel.CutCopyMode = 0

